$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -5
$ws.Range("F7").Value = -6
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -1
$ws.Range("F14").Value = -5
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -2
$ws.Range("F21").Value = -10

$wb.Save()
